$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "kk"
$ws.Range("B7").Value = "kk@gmail.com"
$ws.Range("C7").Value = "kk123"
